$d = $word.ActiveDocument

# 1. Insert a new plain paragraph at the very start of the document containing
#    "<1T>" and move the "_GoBack" bookmark into it.
$startRange = $d.Range(0, 0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>&lt;1T&gt;</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$startRange.InsertXML($newParaXml)

# 2. Collapse the heading paragraph's runs (which were split up by proofErr
#    markers) back into a single run with the same visible text.
$d.Content.Find.Execute("d. has_key( key )", $true, $false, $false, $false, $false, $true, 1, $false, "d. has_key( key )", 2)

# 3. Same collapse for the description paragraph.
$d.Content.Find.Execute("If there is an entry in the dict with the given key , return True, otherwise return False.", $true, $false, $false, $false, $false, $true, 1, $false, "If there is an entry in the dict with the given key , return True, otherwise return False.", 2)

# 4. The old trailing paragraph used to carry the "_GoBack" bookmark; now that
#    the bookmark lives on the new first paragraph, empty it out completely.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$lastPara.Range.InsertXML($emptyParaXml)

Write-Output "done"
